$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 18-34 (team, year, third_place_count, top4_third_place_count)
$data = @(
    @(18, "Egypt",                1990, 27, 0),
    @(23, "Republic of Ireland",  1990,  2, 0),
    @(24, "United States",        1994,  1, 1),
    @(25, "Cameroon",             1994, 12, 12),
    @(26, "South Korea",          1994, 31, 31),
    @(27, "Bulgaria",             1994, 28, 28),
    @(28, "Italy",                1994, 30, 0),
    @(29, "Netherlands",          1994, 25, 0),
    @(30, "Switzerland",          1994, 30, 30),
    @(31, "Russia",               1994, 19, 19),
    @(32, "Norway",               1994,  1, 0),
    @(33, "Belgium",              1994,  6, 0),
    @(34, "Argentina",            1994,  3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Remove now-obsolete trailing rows 35 and 36 (table shrank from 36 to 34 rows)
$ws.Range("A36:D36").EntireRow.Delete()
$ws.Range("A35:D35").EntireRow.Delete()
